# Unity Gantt Chart update
# - Vision Processing Code (GC) task compiled/finished and removed (MATLAB code integrated into Unity)
# - Subsequent rows shift up by one
# - Various plan/actual duration + percent-complete values updated
# - Top "today" marker (H2) and active selection updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the two conditional-formatting ranges that cover the task rows
# (H5:BD26 -> H5:BD25, B27:BD27 -> B26:BD26) before the row shift so the
# grouped cfRules move together.
$ws.Range("H5:BD26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H5:BD25")) | Out-Null
$ws.Range("B27:BD27").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B26:BD26")) | Out-Null

# Remove the "Vision Processing Code (GC)" row entirely (row 14); this also
# shifts rows 15-26 up to 14-25 and drops the now-unused shared string.
$ws.Rows.Item(14).Delete()

# "Today" marker moved from day 11 to day 14
$ws.Range("H2").Value2 = 14

# Row 9 - Acquire VR Setup: actual duration 9 -> 12
$ws.Range("F9").Value2 = 12

# Row 10 - Acquire Vision Software: % complete 0.8 -> 1 (done)
$ws.Range("G10").Value2 = 1

# Row 12 - Skeleton Vision Code: % complete 0.8 -> 1 (done)
$ws.Range("G12").Value2 = 1

# Row 13 - VR Interface: actual duration 7 -> 10, % complete 0.4 -> 0.5
$ws.Range("F13").Value2 = 10
$ws.Range("G13").Value2 = 0.5

# Row 14 - Vision Processing Code (VR): now has actual start/duration and % complete
$ws.Range("E14").Value2 = 10
$ws.Range("F14").Value2 = 5
$ws.Range("G14").Value2 = 0.5

# Row 16 - Vision Processing Code (MATLAB): now has actual start/duration and % complete
$ws.Range("E16").Value2 = 11
$ws.Range("F16").Value2 = 4
$ws.Range("G16").Value2 = 0.4

# Update the active cell selection shown when the sheet is opened
$ws.Range("G15").Select() | Out-Null
